$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove two rows from the existing 15-row detail table (16:30) so it
#     becomes 13 rows (16:28). Shift remaining rows up so the last
#     (specially-bordered) row keeps its closing style and moves from 30 to 28.
$ws.Range("B16:J17").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# --- Update the summary header values ---
$ws.Range("E11").Value = 1268929   # VALOR MORA
$ws.Range("C13").Value = 9         # Cant. Trabajadores
$ws.Range("F13").Value = 11        # Cant. Periodos

# --- Overwrite the detail table (rows 16-28) with the new account-statement data ---
$ws.Cells.Item(16,2).Value = "CC"
$ws.Cells.Item(16,3).Value = "1128056557"
$ws.Cells.Item(16,4).Value = "RAUL EDUARDO LEON GUERRA"
$ws.Cells.Item(16,5).Value = "1706"
$ws.Cells.Item(16,6).Value = 1584
$ws.Cells.Item(16,7).Value = 1264900

$ws.Cells.Item(17,2).Value = "CC"
$ws.Cells.Item(17,3).Value = "9299814"
$ws.Cells.Item(17,4).Value = "DAIRO DE JESUS ARNEDO TORRES"
$ws.Cells.Item(17,5).Value = "1706"
$ws.Cells.Item(17,6).Value = 2864
$ws.Cells.Item(17,7).Value = 2148200

$ws.Cells.Item(18,2).Value = "CC"
$ws.Cells.Item(18,3).Value = "1047410630"
$ws.Cells.Item(18,4).Value = "JOHN FABER LOAIZA FERNANDEZ"
$ws.Cells.Item(18,5).Value = "1706"
$ws.Cells.Item(18,6).Value = 1584
$ws.Cells.Item(18,7).Value = 1734921

$ws.Cells.Item(19,2).Value = "CC"
$ws.Cells.Item(19,3).Value = "71267288"
$ws.Cells.Item(19,4).Value = "JORGE ANDRES GIL RESTREPO"
$ws.Cells.Item(19,5).Value = "1901"
$ws.Cells.Item(19,6).Value = 58447
$ws.Cells.Item(19,7).Value = 2636054

$ws.Cells.Item(20,2).Value = "CC"
$ws.Cells.Item(20,3).Value = "98666591"
$ws.Cells.Item(20,4).Value = "FERNANDO JOSE MONTOYA MORENO"
$ws.Cells.Item(20,5).Value = "1910"
$ws.Cells.Item(20,6).Value = 161808
$ws.Cells.Item(20,7).Value = 5945774

$ws.Cells.Item(21,2).Value = "CC"
$ws.Cells.Item(21,3).Value = "98666591"
$ws.Cells.Item(21,4).Value = "FERNANDO JOSE MONTOYA MORENO"
$ws.Cells.Item(21,5).Value = "1911"
$ws.Cells.Item(21,6).Value = 161808
$ws.Cells.Item(21,7).Value = 5945774

$ws.Cells.Item(22,2).Value = "CC"
$ws.Cells.Item(22,3).Value = "70256194"
$ws.Cells.Item(22,4).Value = "SERGIO ANDRES LONDOÑO CARVAJAL"
$ws.Cells.Item(22,5).Value = "2109"
$ws.Cells.Item(22,6).Value = 7428
$ws.Cells.Item(22,7).Value = 1856933

$ws.Cells.Item(23,2).Value = "CC"
$ws.Cells.Item(23,3).Value = "1050949043"
$ws.Cells.Item(23,4).Value = "GUILLERMO ENRIQUE BAENA RODRIGUEZ"
$ws.Cells.Item(23,5).Value = "2203"
$ws.Cells.Item(23,6).Value = 51917
$ws.Cells.Item(23,7).Value = 1842486

$ws.Cells.Item(24,2).Value = "CE"
$ws.Cells.Item(24,3).Value = "627034"
$ws.Cells.Item(24,4).Value = "JUAN MANUEL NAVARRO ESCOBAR"
$ws.Cells.Item(24,5).Value = "2404"
$ws.Cells.Item(24,6).Value = 398383
$ws.Cells.Item(24,7).Value = 17784960

$ws.Cells.Item(25,2).Value = "CC"
$ws.Cells.Item(25,3).Value = "71360602"
$ws.Cells.Item(25,4).Value = "WILLINGTON GOMEZ OROZCO NIETO"
$ws.Cells.Item(25,5).Value = "2408"
$ws.Cells.Item(25,6).Value = 83110
$ws.Cells.Item(25,7).Value = 3008951

$ws.Cells.Item(26,2).Value = "CC"
$ws.Cells.Item(26,3).Value = "71360602"
$ws.Cells.Item(26,4).Value = "WILLINGTON GOMEZ OROZCO NIETO"
$ws.Cells.Item(26,5).Value = "2409"
$ws.Cells.Item(26,6).Value = 113332
$ws.Cells.Item(26,7).Value = 3008951

$ws.Cells.Item(27,2).Value = "CC"
$ws.Cells.Item(27,3).Value = "71360602"
$ws.Cells.Item(27,4).Value = "WILLINGTON GOMEZ OROZCO NIETO"
$ws.Cells.Item(27,5).Value = "2410"
$ws.Cells.Item(27,6).Value = 113332
$ws.Cells.Item(27,7).Value = 3008951

$ws.Cells.Item(28,2).Value = "CC"
$ws.Cells.Item(28,3).Value = "71360602"
$ws.Cells.Item(28,4).Value = "WILLINGTON GOMEZ OROZCO NIETO"
$ws.Cells.Item(28,5).Value = "2411"
$ws.Cells.Item(28,6).Value = 113332
$ws.Cells.Item(28,7).Value = 3008951
